$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new data row (row 8): time value 0.5 and shared string "c"
$ws.Cells.Item(8, 1).Value = 0.5
$ws.Cells.Item(8, 2).Value = "c"

# Turn off the existing autofilter so it can be reapplied over the
# expanded A1:B8 range (re-invoking AutoFilter on an active filter
# keeps the old range, so drop it first).
$ws.AutoFilterMode = $false

# Reapply the autofilter over the new A1:B8 range with the extra
# filtered value (0.500) alongside the original 0.046 / 0.516 ones.
$rng = $ws.Range("A1:B8")
$rng.AutoFilter(1, @("0.046", "0.500", "0.516"), 7) | Out-Null

# Keep the _xlnm._FilterDatabase defined name in sync with the new
# autofilter range.
for ($i = 1; $i -le $wb.Names.Count; $i++) {
    $n = $wb.Names.Item($i)
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "=Munka1!`$A`$1:`$B`$8"
    }
}

# Move the active selection to C7.
$ws.Range("C7").Select() | Out-Null
